$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row's formatting down into the new row so the new cells
# pick up the same style (s="2") used throughout the header row.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

# Populate the new data row (row 2) with the observation profile details.
$ws.Range("A2").Value = "us-core-observation-adi-documentation"
$ws.Range("B2").Value = "US Core Observation ADI Documentation Profile"
$ws.Range("C2").Value = "null#observation-adi-documentation"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "LOINC#42348-3"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "dateTimeĵ, Periodĵ, Timingĵ, instantĵ"
$ws.Range("H2").Value = "CodeableConceptĵ"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
